# Automatic update of files.
# The source "Avverkningsanmälningar" feed was re-pulled: every record's
# "Förändrad" (C) date advances by one day, and the individual case
# records (rows 4-33) get reshuffled / refreshed against the latest feed
# order, with a couple of cases being replaced by newer ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Every data row's "Förändrad" column (C) moves from 46077 to 46078 ---
for ($r = 2; $r -le 33; $r++) {
    $ws.Cells.Item($r, 3).Value = 46078
}

# --- 2. Row-level refresh: Beteckning (A), Datum (B), Area (G) ---
$data = @(
    @{ Row = 4;  A = "A 389-2023";   B = 44929;              G = 2.5 },
    @{ Row = 6;  A = "A 1782-2024";  B = 45307;              G = 2.7 },
    @{ Row = 8;  A = "A 4481-2024";  B = 45327;              G = 1 },
    @{ Row = 9;  A = "A 18327-2025"; B = 45762;              G = 0.6 },
    @{ Row = 10; A = "A 4493-2024";  B = 45327;              G = 1.8 },
    @{ Row = 11; A = "A 18328-2025"; B = 45762;              G = 1.8 },
    @{ Row = 12; A = "A 10710-2025"; B = 45722;              G = 1.8 },
    @{ Row = 13; A = "A 28260-2023"; B = 45099;              G = 5 },
    @{ Row = 14; A = "A 4822-2023";  B = 44957;              G = 2.2 },
    @{ Row = 15; A = "A 34400-2025"; B = 45846.61351851852;  G = 1.3 },
    @{ Row = 16; A = "A 34401-2025"; B = 45846.6140162037;   G = 2.8 },
    @{ Row = 17; A = "A 1531-2022";  B = 44573;              G = 1.6 },
    @{ Row = 19; A = "A 24-2023";    B = 44928;              G = 0.5 },
    @{ Row = 21; A = "A 4486-2024";  B = 45327;              G = 0.6 },
    @{ Row = 22; A = "A 7727-2026";  B = 46062.50420138889;  G = 1.9 },
    @{ Row = 23; A = "A 7731-2026";  B = 46062.52008101852;  G = 5.9 },
    @{ Row = 24; A = "A 21572-2023"; B = 45063;              G = 1.7 },
    @{ Row = 25; A = "A 635-2023";   B = 44930;              G = 0.5 },
    @{ Row = 26; A = "A 32610-2024"; B = 45513;              G = 0.5 },
    @{ Row = 27; A = "A 18434-2023"; B = 45042;              G = 0.7 },
    @{ Row = 28; A = "A 4487-2024";  B = 45327;              G = 1.9 },
    @{ Row = 29; A = "A 5817-2025";  B = 45694.74113425926;  G = 1.2 },
    @{ Row = 30; A = "A 48974-2023"; B = 45209;              G = 4.5 },
    @{ Row = 31; A = "A 48181-2024"; B = 45589;              G = 0.7 },
    @{ Row = 32; A = "A 53131-2021"; B = 44468;              G = 1.3 },
    @{ Row = 33; A = "A 18332-2025"; B = 45762;              G = 2.5 }
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 1).Value = $item.A
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 7).Value = $item.G
}

# --- 3. Markägare (F) column: the "Kommuner" tag moves from row 30 to row 12 ---
$ws.Range("F12").Value = "Kommuner"
$ws.Range("F30").Value = ""

# --- 4. Rows 4 & 6 also carry HYPERLINK formulas (S,T,V,W,X,Y) that embed
#        the Beteckning text - refresh those to match the new A4/A6 values ---
$linkCols = @(
    @{ Col = "S"; Folder = "artfynd";        Suffix = "artfynd.xlsx" },
    @{ Col = "T"; Folder = "kartor";         Suffix = "karta.png" },
    @{ Col = "V"; Folder = "klagomål";       Suffix = "FSC-klagomål.docx" },
    @{ Col = "W"; Folder = "klagomålsmail";  Suffix = "FSC-klagomål mail.docx" },
    @{ Col = "X"; Folder = "tillsyn";        Suffix = "tillsynsbegäran.docx" },
    @{ Col = "Y"; Folder = "tillsynsmail";   Suffix = "tillsynsbegäran mail.docx" }
)

$linkRows = @(
    @{ Row = 4; A = "A 389-2023" },
    @{ Row = 6; A = "A 1782-2024" }
)

foreach ($lr in $linkRows) {
    foreach ($lc in $linkCols) {
        $url = "https://klasma.github.io/Logging_1278/" + $lc.Folder + "/" + $lr.A + " " + $lc.Suffix
        $formula = '=HYPERLINK("' + $url + '", "' + $lr.A + '")'
        $ws.Range($lc.Col + $lr.Row).Formula = $formula
    }
}
